# Update LDLC prices history:
# Insert a new timestamped price-history column (GV) right before the
# "nom" / "url_produit" columns, shifting them one column to the right
# (nom: GV -> GW, url_produit: GW -> GX).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column 204 = GV. Inserting here pushes the existing GV ("nom") to GW
# and the existing GW ("url_produit") to GX, while preserving styles.
$ws.Columns.Item(204).Insert()

# Header: new snapshot timestamp for the newly inserted column.
$ws.Range("GV1").Value2 = "2026-02-06 11:24:50"

# Data rows: the new column gets a copy of the latest existing price
# snapshot (column GU, now column 203) for rows where that snapshot has
# a numeric price; rows without a recorded price are left blank, same
# as the source column.
for ($r = 2; $r -le 210; $r++) {
    $src = $ws.Cells.Item($r, 203)
    if ($src.Value2 -ne "") {
        $ws.Cells.Item($r, 204).Value2 = $src.Value2
    }
}
